$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

foreach ($row in 2..12) {
    $ws.Cells.Item($row, 3).Value = 45224
}
